$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ativação date update
$ws.Range("B8").Value = "01/01/2023"
$ws.Range("C8").Value = "01/01/2023"

# Objetivos (responsible professor) update
$ws.Range("B10").Value = "7290967 - Emerson Gonçalves de Melo"
$ws.Range("C10").Value = "7290967 - Emerson Gonçalves de Melo"

# Objectives: English description text (new)
$ws.Range("B11").Value = "Provide means for the student to acquire knowledge and practice in the use of modern computational tools applied to physical problems."
$ws.Range("C11").Value = "Provide means for the student to acquire knowledge and practice in the use of modern computational tools applied to physical problems."

# Programa resumido (mirrors the Ativação date, per source data)
$ws.Range("B13").Value = "01/01/2023"
$ws.Range("C13").Value = "01/01/2023"

# Short syllabus: English text (new)
$ws.Range("B14").Value = "Numerical simulation in deterministic and stochastic systems. Monte Carlo methods. Random walks. Fractals. Introduction to spectral analysis by Fourier transforms. Review of solutions of ordinary and partial differential equations. Numerical solution of partial differential equations."
$ws.Range("C14").Value = "Numerical simulation in deterministic and stochastic systems. Monte Carlo methods. Random walks. Fractals. Introduction to spectral analysis by Fourier transforms. Review of solutions of ordinary and partial differential equations. Numerical solution of partial differential equations."

# Programa (mirrors the responsible professor, per source data)
$ws.Range("B15").Value = "7290967 - Emerson Gonçalves de Melo"
$ws.Range("C15").Value = "7290967 - Emerson Gonçalves de Melo"

# Syllabus: bullet-point English program text (new)
$ws.Range("B16").Value = "• Numerical simulation in deterministic systems. • Numerical simulations in stochastic systems. • Pseudo-random numbers. • Random walks, diffusion and percolation.  • Fractals. • Introduction to spectral analysis by Fourier transforms. • Fast Fourier Transform and applications. • Identification of frequencies and normal modes. • Signal detection and treatment • Image processing. • Numerical solution methods for partial differential equations. • Finite difference method. • Finite Element Method."
$ws.Range("C16").Value = "• Numerical simulation in deterministic systems. • Numerical simulations in stochastic systems. • Pseudo-random numbers. • Random walks, diffusion and percolation.  • Fractals. • Introduction to spectral analysis by Fourier transforms. • Fast Fourier Transform and applications. • Identification of frequencies and normal modes. • Signal detection and treatment • Image processing. • Numerical solution methods for partial differential equations. • Finite difference method. • Finite Element Method."

# Método: responsible professor (replacing the removed co-professor)
$ws.Range("B18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
